# Convert an "RRGGBB" hex string into the packed integer PowerPoint's
# RGB color properties expect (little-endian / BGR-ordered: R + G*256 + B*65536).
function Get-RGBValue($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table on slide 16: switch its table style to the built-in
#        "Medium Style 2 - Accent 1" style. ---------------------------------
$s16 = $p.Slides.Item(16)
$tableShape = $s16.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{61ECFCA6-67BE-4FAF-BDE2-B22B47BDE962}")

# --- 2. Design/theme colour scheme: swap the "Integral" theme palette for
#        the plain "Office" palette. -----------------------------------------
$officeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$themeColors = $s16.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = Get-RGBValue($officeColors[$i - 1])
}
